# The deck currently carries the "Integral" theme on its slide master
# (ppt/theme/theme1.xml) and the default "Office Theme" on its notes
# master (ppt/theme/theme2.xml). The edit re-colors the slide master's
# theme with the standard Office color palette (matching what used to
# live on the notes-master theme), i.e. the deck's main theme colors
# change from the green/olive "Integral" palette to the blue/orange
# default "Office" palette.
#
# PowerPoint's Design/Theme COM surface does not offer a bulk "apply a
# named color scheme" call, so each of the twelve theme colors is set
# individually via ColorScheme.Colors(index).RGB, in the standard
# ppColorSchemeIndex order:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2,
#   5 accent1, 6 accent2, 7 accent3, 8 accent4, 9 accent5, 10 accent6,
#   11 hlink, 12 folHlink
# RGB is a VBA-style BGR long (&H00BBGGRR), so each hex literal below
# is the target RRGGBB value byte-reversed.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.ColorScheme

$scheme.Colors(1).RGB  = 0x000000   # dk1      -> 000000
$scheme.Colors(2).RGB  = 0xFFFFFF   # lt1      -> FFFFFF
$scheme.Colors(3).RGB  = 0x6A5444   # dk2      -> 44546A
$scheme.Colors(4).RGB  = 0xE6E6E7   # lt2      -> E7E6E6
$scheme.Colors(5).RGB  = 0xD59B5B   # accent1  -> 5B9BD5
$scheme.Colors(6).RGB  = 0x317DED   # accent2  -> ED7D31
$scheme.Colors(7).RGB  = 0xA5A5A5   # accent3  -> A5A5A5
$scheme.Colors(8).RGB  = 0x00C0FF   # accent4  -> FFC000
$scheme.Colors(9).RGB  = 0xC47244   # accent5  -> 4472C4
$scheme.Colors(10).RGB = 0x47AD70   # accent6  -> 70AD47
$scheme.Colors(11).RGB = 0xC16305   # hlink    -> 0563C1
$scheme.Colors(12).RGB = 0x724F95   # folHlink -> 954F72
